$d = $word.ActiveDocument

# 1. Update activation date (Ativação: 01/01/2022 -> 01/01/2024)
$d.Content.Find.Execute(
    "Ativação: 01/01/2022", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2024", 2
) | Out-Null

# 2. Add a new co-instructor on its own line after Gilberto Carvalho Coelho.
#    We insert a manual line break + the new name right after the existing
#    text, then nudge formatting on the inserted text so the engine keeps
#    it in its own <w:r> run (matching the two-run shape of the diff)
#    instead of silently merging it back into the previous run.
$gilbertoPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*5009972 - Gilberto Carvalho Coelho*") {
        $gilbertoPara = $p
        break
    }
}
$newName = "984972 - Hugo Ricardo Zschommler Sandim"
$pRange = $gilbertoPara.Range
$textRange = $d.Range($pRange.Start, $pRange.End - 1)
$insertionPoint = $textRange.End
$textRange.InsertAfter([char]11 + $newName)
$newNameRange = $d.Range($insertionPoint + 1, $insertionPoint + 1 + $newName.Length)
$newNameRange.Bold = 1
$newNameRange.Bold = 0

# 3. Update "Método" paragraph text
$d.Content.Find.Execute(
    "O curso será ministrado na forma de aulas expositivas. Estão previstas visitas a empresas de fundição para consolidação dos conceitos apresentados nas aulas expositivas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "O curso será ministrado na forma de aulas expositivas e visitas a empresas de fundição. Questionários e listas de exercícios serão elaborados para serem respondidos individualmente ou em grupo. Avaliações escritas serão realizadas para resolução individual.",
    2
) | Out-Null

# 4. Update "Critério" paragraph text
$d.Content.Find.Execute(
    "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). O critério para a nota final é: NF=(P1+P2)/2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As avaliações individuais, a participação nas resoluções dos exercícios e repostas aos questionários serão agrupadas em duas notas (N1 e N2) que comporão a nota final (NF). O critério para cálculo da nota final é: NF = (N1+ N2)/2Serão aprovados os alunos com NF ≥ 5,0Serão reprovados os alunos com NF < 3,0",
    2
) | Out-Null

# 5. Update "Norma de recuperação" paragraph text
$d.Content.Find.Execute(
    "Será aplicada uma prova de recuperação cuja nota comporá média aritmética com a nota final NF.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Será aplicada recuperação para os alunos que obtiverem NF entre 3,0 e 4,9. A nota pós recuperação será calculada pela média aritmética com a nota final NF.",
    2
) | Out-Null
